# Applies the "Added entry for Bryan Alcala" commit.
#
# Strategy: several existing runs of text get split into multiple runs,
# with w:proofErr (grammar-check) markers wrapped around certain words.
# We locate each original run of text with Find, then replace the precise
# sub-range using Range.InsertXML with hand-built OOXML fragments that
# reproduce the exact <w:r>/<w:proofErr> structure from the target diff.
# Finally, four new paragraphs (Bryan Alcala's entry) are inserted before
# the document's final (blank) paragraph.

$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>'

function Run($text, $preserve) {
    if ($preserve) {
        return "<w:r>$rPr<w:t xml:space=`"preserve`">$text</w:t></w:r>"
    } else {
        return "<w:r>$rPr<w:t>$text</w:t></w:r>"
    }
}

function GramRun($text) {
    return "<w:proofErr w:type=`"gramStart`"/>" + (Run $text) + "<w:proofErr w:type=`"gramEnd`"/>"
}

function Wrap($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function ReplaceFound([string]$searchText, [string]$replacementXml) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $searchText"
        return
    }
    # IMPORTANT: build a *fresh* Range object from the found Start/End
    # rather than reusing $r (which Find.Execute mutated in place).
    # Re-using the mutated range causes InsertXML to insert beside the
    # text instead of replacing it.
    $target = $d.Range($r.Start, $r.End)
    $xml = Wrap("<w:body><w:p>$replacementXml</w:p></w:body>")
    $target.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1. "add today's date and your name" (bullet list item)
# ---------------------------------------------------------------------
$xml1 = (Run "add today's date and your " $true) + (GramRun "name")
ReplaceFound "add today's date and your name" $xml1

# ---------------------------------------------------------------------
# 2. "Add a line with today's date and your name" (bullet list item)
# ---------------------------------------------------------------------
$xml2 = (Run "Add a line with today's date and your " $true) + (GramRun "name")
ReplaceFound "Add a line with today's date and your name" $xml2

# ---------------------------------------------------------------------
# 3. David Singletary bio: "...I love technology and coding, my favorite..."
# ---------------------------------------------------------------------
$xml3 = (Run "I am a retired software engineer who loves books, movies, and television shows (especially classic horror/sci-fi). I love technology and " $true) `
      + (GramRun "coding,") `
      + (Run " my favorite languages are Java and C++ but I also teach R and Python in our Data Science program." $true)
ReplaceFound "I am a retired software engineer who loves books, movies, and television shows (especially classic horror/sci-fi). I love technology and coding, my favorite languages are Java and C++ but I also teach R and Python in our Data Science program." $xml3

# ---------------------------------------------------------------------
# 4. Lilly Nguyen bio, part 1: "Hello, my name is Lilly Nguyen and I am born..."
# ---------------------------------------------------------------------
$xml4 = (Run "Hello, my name is Lilly " $true) `
      + (GramRun "Nguyen") `
      + (Run " and I " $true) `
      + (GramRun "am") `
      + (Run " born and raised in Jacksonville FL. I enjoy playing video games and fishing in my spare time. I enjoy playing video games such as " $true)
ReplaceFound "Hello, my name is Lilly Nguyen and I am born and raised in Jacksonville FL. I enjoy playing video games and fishing in my spare time. I enjoy playing video games such as " $xml4

# ---------------------------------------------------------------------
# 5. Lilly Nguyen bio, part 2: ", TFT, and BTD6. ... eventually fish an alligator."
# ---------------------------------------------------------------------
$xml5 = (Run ", TFT, and BTD6. As for fishing, I mostly do saltwater fishing and have found the magical place called the Jax Pier. I hope to eventually " $true) `
      + (GramRun "fish") `
      + (Run " an alligator. " $true)
ReplaceFound ", TFT, and BTD6. As for fishing, I mostly do saltwater fishing and have found the magical place called the Jax Pier. I hope to eventually fish an alligator. " $xml5

# ---------------------------------------------------------------------
# 6. Append Bryan Alcala's entry before the final (blank) paragraph.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
$jasonBioPara = $paras.Item($count - 1)
$insertPoint = $jasonBioPara.Range.End
$ins = $d.Range($insertPoint, $insertPoint)

$blankPara = "<w:p><w:pPr>$rPr</w:pPr></w:p>"

$dateLine = "<w:p><w:pPr>$rPr</w:pPr>" + (Run "9/5/2023 Bryan Alcala") + "</w:p>"

$bioXml = (Run "Hi everyone, " $true) `
        + (GramRun "My") `
        + (Run " name is Bryan. I am 27 years old and love working with technology. I have an " $true) `
        + (GramRun "Associates in Information Technology") `
        + (Run " and am working my way through my Bachelors. I am excited to start this class and learn with you all!" $true)
$bioLine = "<w:p><w:pPr>$rPr</w:pPr>" + $bioXml + "</w:p>"

$newParas = $blankPara + $dateLine + $bioLine + $blankPara

$finalXml = Wrap("<w:body>$newParas</w:body>")
$ins.InsertXML($finalXml)

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
